# "Updating dataset so I don't have that nonsense with the too-short song"
#
# The experimenter reran the 120g.start / 120g.minute / 120g.mix group with a
# changed dataset. The three 120g.start samples (12gS.0.0, 12gS.0.1, 12gS.0.2 —
# rows 115-117) are no longer valid against the new dataset, so their recorded
# percentage (column C) and source file (column H) are pulled out, and the
# original numbers are re-homed under a brand new group "120g.start0" in three
# freshly appended rows (124-126), using new sample labels 12gS0.0.0/.1/.2.
# A comment on B124 documents why.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log")

# --- 1. Pull the stale results out of rows 115:117 -------------------------
# (group 120g.start no longer matches the new dataset for these 3 runs)
$ws.Range("C115:C117").ClearContents()
$ws.Range("H115:H117").ClearContents()

# --- 2. Re-home that same data under the new "120g.start0" group ----------
$ws.Range("A124").Value = "12gS0.0.0"
$ws.Range("B124").Value = "120g.start0"
$ws.Range("C124").Value = 0.266666669646898
$ws.Range("D124").Value = 0.33
$ws.Range("E124").Formula = "=C124/D124"
$ws.Range("F124").Formula = "=AVERAGEIF(B:B,B124,E:E)"
$ws.Range("G124").Value = 42572
$ws.Range("H124").Value = "09.41.txt"

$ws.Range("A125").Value = "12gS0.0.1"
$ws.Range("B125").Value = "120g.start0"
$ws.Range("C125").Value = 0.26666667064030902
$ws.Range("D125").Value = 0.33
$ws.Range("E125").Formula = "=C125/D125"
$ws.Range("F125").Formula = "=AVERAGEIF(B:B,B125,E:E)"
$ws.Range("G125").Value = 42572
$ws.Range("H125").Value = "09.41.txt"

$ws.Range("A126").Value = "12gS0.0.2"
$ws.Range("B126").Value = "120g.start0"
$ws.Range("C126").Value = 0.56666667858759501
$ws.Range("D126").Value = 0.33
$ws.Range("E126").Formula = "=C126/D126"
$ws.Range("F126").Formula = "=AVERAGEIF(B:B,B126,E:E)"
$ws.Range("G126").Value = 42572
$ws.Range("H126").Value = "09.41.txt"

# Dates in G/H should keep the same date/"text-as-date" display as the rest
# of the column.
$ws.Range("G124:G126").NumberFormat = $ws.Range("G114").NumberFormat
$ws.Range("H124:H126").NumberFormat = $ws.Range("H114").NumberFormat

# --- 3. Explain the re-homing with a cell comment --------------------------
$note = "20g.start0: originally was 120g.start, but threw out these results from that group because I changed the dataset."
$cmt = $ws.Range("B124").AddComment($note)

# --- 4. Match the selection left behind by the edit -------------------------
$ws.Range("G123").Select()
